$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.715.65'
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").Value = '3.500.81'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''593.62'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '''169.04'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.591'
$ws.Range("E8").Value = '  +5.22%  '
$ws.Range("E9").Value = '  +8.12%  '
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").Value = '4.105.89'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  +2.64%  '
$ws.Range("D15").Value = '''0.0000181'
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").Value = '66.734.62'
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("D17").Value = '3.502.09'
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("D19").Value = '''14.05'
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").Value = '''395.99'
$ws.Range("E20").Value = '  +3.12%  '
$ws.Range("D21").Value = '''7.98'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '''73.49'
$ws.Range("E22").Value = '  +2.25%  '
$ws.Range("D24").Value = '''0.537'
$ws.Range("E24").Value = '  +3.01%  '
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("D26").Value = '''10.17'
$ws.Range("E26").Value = '  +1.38%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  -2.36%  '
$ws.Range("D29").Value = '''6.30'
$ws.Range("E29").Value = '  -0.52%  '
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("D32").Value = '''23.83'
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").Value = '''1.62'
$ws.Range("E34").Value = '  +6.21%  '
$ws.Range("D35").Value = '''162.62'
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("D36").Value = '''0.897'
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  +2.53%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''6.83'
$ws.Range("E38").Value = '  +2.77%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").Value = '''4.68'
$ws.Range("E39").Value = '  +4.79%  '
$ws.Range("D40").Value = '''0.0745'
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("D41").Value = '''26.54'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("D42").Value = '''27.03'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").Value = '2.796.57'
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''2.59'
$ws.Range("E44").Value = '  +2.96%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''42.90'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("E46").Value = '  +1.03%  '
$ws.Range("D47").Value = '''342.53'
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("D48").Value = '''1.10'
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("E49").Value = '  +4.54%  '
$ws.Range("D50").Value = '''0.859'
$ws.Range("E50").Value = '  +2.01%  '
$ws.Range("E51").Value = '  +1.82%  '
